{"js": "const body = context.document.body;\nconst pairs = [\n  [\"58+6=64\", \"32+56=88\"],\n  [\"43-5=38\", \"19+2=21\"],\n  [\"89+2=91\", \"82-28=54\"],\n  [\"71-12=59\", \"11+79=90\"],\n  [\"0+6=6\", \"51-5=46\"],\n  [\"64-55=9\", \"78-76=2\"],\n  [\"70-64=6\", \"40-40=0\"],\n  [\"8+52=60\", \"43+35=78\"],\n  [\"50+37=87\", \"12+86=98\"],\n  [\"11+55=66\", \"90-69=21\"],\n  [\"59-23=36\", \"13+51=64\"],\n  [\"95+4=99\", \"78-11=67\"],\n  [\"81-27=54\", \"31-21=10\"],\n  [\"61-41=20\", \"72-10=62\"],\n  [\"17-4=13\", \"60+21=81\"],\n  [\"83-15=68\", \"87-71=16\"],\n  [\"65-17=48\", \"66-39=27\"],\n  [\"38+54=92\", \"45-6=39\"],\n  [\"56-19=37\", \"37-26=11\"],\n  [\"33+36=69\", \"77-50=27\"],\n  [\"33+60=93\", \"76-22=54\"],\n  [\"4+47=51\", \"79-74=5\"],\n  [\"75-33=42\", \"12-4=8\"],\n  [\"49+6=55\", \"16+36=52\"],\n  [\"90-65=25\", \"52+38=90\"],\n  [\"43+3=46\", \"67-27=40\"],\n  [\"80-60=20\", \"80-72=8\"],\n  [\"25-16=9\", \"93-92=1\"],\n  [\"0+76=76\", \"87-10=77\"],\n  [\"11+19=30\", \"22+24=46\"],\n  [\"41-38=3\", \"1+65=66\"],\n  [\"77-55=22\", \"7+58=65\"],\n  [\"31+67=98\", \"57+28=85\"],\n  [\"81-74=7\", \"73-4=69\"],\n  [\"69-51=18\", \"53+39=92\"],\n  [\"24+73=97\", \"37-23=14\"],\n  [\"67-52=15\", \"31+16=47\"],\n  [\"25-6=19\", \"98-61=37\"],\n  [\"14-8=6\", \"86-82=4\"],\n  [\"60-31=29\", \"29+13=42\"],\n  [\"62-14=48\", \"16+31=47\"],\n  [\"92-46=46\", \"74-68=6\"],\n  [\"5+79=84\", \"11-2=9\"],\n  [\"83-38=45\", \"66+14=80\"],\n  [\"94-21=73\", \"36+30=66\"],\n  [\"52+17=69\", \"97-23=74\"],\n  [\"99-4=95\", \"3+91=94\"],\n  [\"84-71=13\", \"45-29=16\"],\n  [\"27+51=78\", \"59-27=32\"],\n  [\"34+22=56\", \"17+35=52\"],\n  [\"23-11=12\", \"59+37=96\"],\n  [\"90+5=95\", \"43-13=30\"],\n  [\"59+1=60\", \"95-23=72\"],\n  [\"61+34=95\", \"52+12=64\"],\n  [\"58-35=23\", \"6+47=53\"],\n  [\"11+33=44\", \"54-4=50\"],\n  [\"9+65=74\", \"94+5=99\"],\n  [\"66-44=22\", \"30+37=67\"],\n  [\"73+25=98\", \"96-15=81\"],\n  [\"98-93=5\", \"85-4=81\"],\n  [\"33+1=34\", \"99-26=73\"],\n  [\"95-29=66\", \"66+3=69\"],\n  [\"57-7=50\", \"63-37=26\"],\n  [\"69-63=6\", \"72-17=55\"],\n  [\"95-85=10\", \"27+31=58\"],\n  [\"42-0=42\", \"84-44=40\"],\n  [\"96-62=34\", \"74-14=60\"],\n  [\"36+50=86\", \"56+6=62\"],\n  [\"32+35=67\", \"2+5=7\"],\n  [\"4+31=35\", \"88-84=4\"],\n  [\"38-17=21\", \"81-45=36\"],\n  [\"6+86=92\", \"73+6=79\"],\n  [\"69-20=49\", \"35+57=92\"],\n  [\"83-58=25\", \"72+19=91\"],\n  [\"46+33=79\", \"90+4=94\"],\n  [\"41+52=93\", \"46-45=1\"],\n  [\"7+25=32\", \"82-45=37\"],\n  [\"92-73=19\", \"98-44=54\"],\n  [\"91-19=72\", \"99-83=16\"],\n  [\"23+20=43\", \"76-28=48\"],\n  [\"34-7=27\", \"54+16=70\"],\n  [\"96-32=64\", \"81-31=50\"],\n  [\"69-0=69\", \"66-62=4\"],\n  [\"93-20=73\", \"32-9=23\"],\n  [\"62+34=96\", \"67+8=75\"],\n  [\"41+15=56\", \"49-27=22\"],\n  [\"47-13=34\", \"58+13=71\"],\n  [\"76-4=72\", \"95-33=62\"],\n  [\"12+62=74\", \"45-8=37\"],\n  [\"71+24=95\", \"64+1=65\"],\n  [\"9+51=60\", \"37-1=36\"],\n  [\"81-54=27\", \"88-6=82\"],\n  [\"12+13=25\", \"22+29=51\"],\n  [\"19-7=12\", \"55+44=99\"],\n  [\"38+21=59\", \"10-1=9\"],\n  [\"29+9=38\", \"87-2=85\"],\n  [\"86-64=22\", \"81-68=13\"],\n  [\"26+50=76\", \"28+54=82\"],\n  [\"21+44=65\", \"90-5=85\"],\n  [\"36-15=21\", \"4+21=25\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"58+6=64\", \"32+56=88\")\n  ,@(\"43-5=38\", \"19+2=21\")\n  ,@(\"89+2=91\", \"82-28=54\")\n  ,@(\"71-12=59\", \"11+79=90\")\n  ,@(\"0+6=6\", \"51-5=46\")\n  ,@(\"64-55=9\", \"78-76=2\")\n  ,@(\"70-64=6\", \"40-40=0\")\n  ,@(\"8+52=60\", \"43+35=78\")\n  ,@(\"50+37=87\", \"12+86=98\")\n  ,@(\"11+55=66\", \"90-69=21\")\n  ,@(\"59-23=36\", \"13+51=64\")\n  ,@(\"95+4=99\", \"78-11=67\")\n  ,@(\"81-27=54\", \"31-21=10\")\n  ,@(\"61-41=20\", \"72-10=62\")\n  ,@(\"17-4=13\", \"60+21=81\")\n  ,@(\"83-15=68\", \"87-71=16\")\n  ,@(\"65-17=48\", \"66-39=27\")\n  ,@(\"38+54=92\", \"45-6=39\")\n  ,@(\"56-19=37\", \"37-26=11\")\n  ,@(\"33+36=69\", \"77-50=27\")\n  ,@(\"33+60=93\", \"76-22=54\")\n  ,@(\"4+47=51\", \"79-74=5\")\n  ,@(\"75-33=42\", \"12-4=8\")\n  ,@(\"49+6=55\", \"16+36=52\")\n  ,@(\"90-65=25\", \"52+38=90\")\n  ,@(\"43+3=46\", \"67-27=40\")\n  ,@(\"80-60=20\", \"80-72=8\")\n  ,@(\"25-16=9\", \"93-92=1\")\n  ,@(\"0+76=76\", \"87-10=77\")\n  ,@(\"11+19=30\", \"22+24=46\")\n  ,@(\"41-38=3\", \"1+65=66\")\n  ,@(\"77-55=22\", \"7+58=65\")\n  ,@(\"31+67=98\", \"57+28=85\")\n  ,@(\"81-74=7\", \"73-4=69\")\n  ,@(\"69-51=18\", \"53+39=92\")\n  ,@(\"24+73=97\", \"37-23=14\")\n  ,@(\"67-52=15\", \"31+16=47\")\n  ,@(\"25-6=19\", \"98-61=37\")\n  ,@(\"14-8=6\", \"86-82=4\")\n  ,@(\"60-31=29\", \"29+13=42\")\n  ,@(\"62-14=48\", \"16+31=47\")\n  ,@(\"92-46=46\", \"74-68=6\")\n  ,@(\"5+79=84\", \"11-2=9\")\n  ,@(\"83-38=45\", \"66+14=80\")\n  ,@(\"94-21=73\", \"36+30=66\")\n  ,@(\"52+17=69\", \"97-23=74\")\n  ,@(\"99-4=95\", \"3+91=94\")\n  ,@(\"84-71=13\", \"45-29=16\")\n  ,@(\"27+51=78\", \"59-27=32\")\n  ,@(\"34+22=56\", \"17+35=52\")\n  ,@(\"23-11=12\", \"59+37=96\")\n  ,@(\"90+5=95\", \"43-13=30\")\n  ,@(\"59+1=60\", \"95-23=72\")\n  ,@(\"61+34=95\", \"52+12=64\")\n  ,@(\"58-35=23\", \"6+47=53\")\n  ,@(\"11+33=44\", \"54-4=50\")\n  ,@(\"9+65=74\", \"94+5=99\")\n  ,@(\"66-44=22\", \"30+37=67\")\n  ,@(\"73+25=98\", \"96-15=81\")\n  ,@(\"98-93=5\", \"85-4=81\")\n  ,@(\"33+1=34\", \"99-26=73\")\n  ,@(\"95-29=66\", \"66+3=69\")\n  ,@(\"57-7=50\", \"63-37=26\")\n  ,@(\"69-63=6\", \"72-17=55\")\n  ,@(\"95-85=10\", \"27+31=58\")\n  ,@(\"42-0=42\", \"84-44=40\")\n  ,@(\"96-62=34\", \"74-14=60\")\n  ,@(\"36+50=86\", \"56+6=62\")\n  ,@(\"32+35=67\", \"2+5=7\")\n  ,@(\"4+31=35\", \"88-84=4\")\n  ,@(\"38-17=21\", \"81-45=36\")\n  ,@(\"6+86=92\", \"73+6=79\")\n  ,@(\"69-20=49\", \"35+57=92\")\n  ,@(\"83-58=25\", \"72+19=91\")\n  ,@(\"46+33=79\", \"90+4=94\")\n  ,@(\"41+52=93\", \"46-45=1\")\n  ,@(\"7+25=32\", \"82-45=37\")\n  ,@(\"92-73=19\", \"98-44=54\")\n  ,@(\"91-19=72\", \"99-83=16\")\n  ,@(\"23+20=43\", \"76-28=48\")\n  ,@(\"34-7=27\", \"54+16=70\")\n  ,@(\"96-32=64\", \"81-31=50\")\n  ,@(\"69-0=69\", \"66-62=4\")\n  ,@(\"93-20=73\", \"32-9=23\")\n  ,@(\"62+34=96\", \"67+8=75\")\n  ,@(\"41+15=56\", \"49-27=22\")\n  ,@(\"47-13=34\", \"58+13=71\")\n  ,@(\"76-4=72\", \"95-33=62\")\n  ,@(\"12+62=74\", \"45-8=37\")\n  ,@(\"71+24=95\", \"64+1=65\")\n  ,@(\"9+51=60\", \"37-1=36\")\n  ,@(\"81-54=27\", \"88-6=82\")\n  ,@(\"12+13=25\", \"22+29=51\")\n  ,@(\"19-7=12\", \"55+44=99\")\n  ,@(\"38+21=59\", \"10-1=9\")\n  ,@(\"29+9=38\", \"87-2=85\")\n  ,@(\"86-64=22\", \"81-68=13\")\n  ,@(\"26+50=76\", \"28+54=82\")\n  ,@(\"21+44=65\", \"90-5=85\")\n  ,@(\"36-15=21\", \"4+21=25\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $rng = $d.Content\n  $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"No match found for: $oldText\"\n  }\n}"}
